{"js": "// \"Align title and dateline flush-left\"\n//\n// The document's Title paragraph (\"Stockholder Notice Under Delaware\n// General Corporation Law 228(e)\") uses the built-in \"Title\" style, and\n// the dateline paragraph (\"January 1, 2016\") uses the built-in \"Date\"\n// style. Both styles currently center their paragraphs; switch their\n// paragraph justification to left so the title and dateline sit flush\n// against the left margin instead of being centered.\n//\n// Re-point the style definitions themselves (rather than the individual\n// paragraphs) since the alignment lives in the style, exactly like a user\n// choosing Home > Paragraph > Align Left while the Title/Date style is\n// selected in the Styles pane (Modify Style), which is a style-level edit.\n\nconst styles = context.document.getStyles();\nconst titleStyle = styles.getByNameOrNullObject(\"Title\");\nconst dateStyle = styles.getByNameOrNullObject(\"Date\");\ntitleStyle.load(\"isNullObject\");\ndateStyle.load(\"isNullObject\");\nawait context.sync();\n\nif (!titleStyle.isNullObject) {\n  titleStyle.paragraphFormat.alignment = Word.Alignment.left;\n}\nif (!dateStyle.isNullObject) {\n  dateStyle.paragraphFormat.alignment = Word.Alignment.left;\n}\n\nawait context.sync();\n", "ps1": "# \"Align title and dateline flush-left\"\n#\n# The document's Title paragraph (\"Stockholder Notice Under Delaware\n# General Corporation Law 228(e)\") uses the built-in \"Title\" style, and\n# the dateline paragraph (\"January 1, 2016\") uses the built-in \"Date\"\n# style. Both styles currently center their paragraphs; switch their\n# paragraph justification to left so the title and dateline sit flush\n# against the left margin instead of being centered.\n#\n# Re-point the style definitions themselves (rather than the individual\n# paragraphs) since the alignment lives in the style, exactly like a user\n# choosing Home > Paragraph > Align Left while the Title/Date style is\n# selected (Modify Style), which is a style-level edit.\n\n$d = $word.ActiveDocument\n\n$wdAlignParagraphLeft = 0\n\ntry {\n  $titleStyle = $d.Styles(\"Title\")\n  $titleStyle.ParagraphFormat.Alignment = $wdAlignParagraphLeft\n} catch {\n}\n\ntry {\n  $dateStyle = $d.Styles(\"Date\")\n  $dateStyle.ParagraphFormat.Alignment = $wdAlignParagraphLeft\n} catch {\n}\n"}
